$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.752.08'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '3.488.75'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.34'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.93'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.25'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('D12').Value = '4.092.40'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.79'
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('D15').Value = '66.788.34'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = '3.498.04'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '392.35'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.63'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000120'
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.19'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.22'
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -3.46%  '
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.38'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.94'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.67'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.30'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0741'
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.15'
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('D44').Value = '2.802.12'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.57'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.55'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '336.99'
$ws.Range('E48').Value = '  -3.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.33'
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('E51').Value = '  -1.37%  '
